$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and E retain their original text formatting so that
# decimal-looking strings (e.g. "0.9992") are not reinterpreted as numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.284.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.843.64'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.06%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.12'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6721'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.19%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07438'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.50%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.89'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07720'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.838.00'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.010'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.94%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6713'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '86.02'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.164'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.220.52'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008322'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '228.71'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.19%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.163'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '160.93'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.706'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.61%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1404'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -3.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.03'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.29%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.162'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.067'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.194'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05288'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.874'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7508'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.97%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.33%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.679'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.320.98'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01807'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.728'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.18%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9201'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.971'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.08341'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +13.11%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.15'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.970.07'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5163'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.65%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '63.83'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.83%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000121'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.140'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05945'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.06%  '
